$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BO1").Value = 0.829539891452441
$ws.Range("BP1").Value = 0.82007215063909134
$ws.Range("A2").Value = 0.80259373283927027
$ws.Range("C2").Value = 0.93071985496064302
$ws.Range("BP2").Value = 0.97385068678740483
$ws.Range("A3").Value = 0.69355162180368857
$ws.Range("D3").Value = 0.98852961852102428
$ws.Range("Z3").Value = 0.82391785696168707
$ws.Range("B4").Value = 0.80492452796392255
$ws.Range("AI4").Value = 0.75623441669188773
$ws.Range("C5").Value = 0.98148201126485102
$ws.Range("F5").Value = 0.65357974899882831
$ws.Range("D6").Value = 0.696592367604149
$ws.Range("F7").Value = 0.60415014627267183
$ws.Range("AU7").Value = 0.81832694000908579
$ws.Range("G8").Value = 0.72358225242784113
$ws.Range("I8").Value = 0.65781110912536001
$ws.Range("J8").Value = 0.85606342206580055
$ws.Range("G9").Value = 0.99507244005257367
$ws.Range("J9").Value = 0.96030110100067911
$ws.Range("K9").Value = 0.95451697391314361
$ws.Range("R9").Value = 0.99020287016935982
$ws.Range("AF9").Value = 0.94425384632779408
$ws.Range("L10").Value = 0.84495696841767931
$ws.Range("J11").Value = 0.59810723282492173
$ws.Range("AY11").Value = 0.93812181784764892
$ws.Range("K12").Value = 0.86062237429283395
$ws.Range("N12").Value = 0.68380533002232702
$ws.Range("Z12").Value = 0.89680406480336805
$ws.Range("M14").Value = 0.946499108099216
$ws.Range("P14").Value = 0.67819448867337617
$ws.Range("T14").Value = 0.98873415200878512
$ws.Range("M15").Value = 0.74840995157321866
$ws.Range("N15").Value = 0.84759506097686521
$ws.Range("Q15").Value = 0.86597837169010328
$ws.Range("BE17").Value = 0.9393294330387163
$ws.Range("P18").Value = 0.92957678853152537
$ws.Range("BG18").Value = 0.91894999563539792
$ws.Range("Q19").Value = 0.99689920105526919
$ws.Range("R19").Value = 0.95767668104906345
$ws.Range("U19").Value = 0.85311497721446872
$ws.Range("S20").Value = 0.90187095963262087
$ws.Range("T21").Value = 0.88663132798935029
$ws.Range("AP21").Value = 0.56891781421045229
$ws.Range("U22").Value = 0.8636615500167647
$ws.Range("X22").Value = 0.65798446133782029
$ws.Range("AF22").Value = 0.82447853089763967
$ws.Range("AM22").Value = 0.9143429622965662
$ws.Range("U23").Value = 0.95511581071439033
$ws.Range("W24").Value = 0.99489612008245487
$ws.Range("Y24").Value = 0.9484819974786074
$ws.Range("Z24").Value = 0.79542623672467938
$ws.Range("W25").Value = 0.79698105988083268
$ws.Range("BG25").Value = 0.77093728812163742
$ws.Range("AX26").Value = 0.96930723592938195
$ws.Range("AC27").Value = 0.79245131980508798
$ws.Range("BD27").Value = 0.74770457770775467
$ws.Range("Z28").Value = 0.98913755805847114
$ws.Range("AC28").Value = 0.97709654936972679
$ws.Range("AD28").Value = 0.74939096413429707
$ws.Range("E29").Value = 0.87706050415019676
$ws.Range("AS29").Value = 0.81495842467671165
$ws.Range("AF30").Value = 0.88581646743325748
$ws.Range("AD31").Value = 0.98541560758090285
$ws.Range("AF31").Value = 0.83650767673215287
$ws.Range("AX31").Value = 0.85809330897567193
$ws.Range("AO32").Value = 0.63361570587124438
$ws.Range("AH33").Value = 0.85893276212978642
$ws.Range("AI33").Value = 0.80213773010997746
$ws.Range("BD33").Value = 0.96433013033600845
$ws.Range("AJ34").Value = 0.96908537917190651
$ws.Range("AH35").Value = 0.59804738369992472
$ws.Range("AJ35").Value = 0.7560095305572001
$ws.Range("AK35").Value = 0.81776173449418277
$ws.Range("AK36").Value = 0.96083967433405504
$ws.Range("AL36").Value = 0.82736360102386386
$ws.Range("AL37").Value = 0.77909097264914928
$ws.Range("AO37").Value = 0.89365137467350308
$ws.Range("AM38").Value = 0.96519269540330666
$ws.Range("AK39").Value = 0.75192975494155856
$ws.Range("Y40").Value = 0.55970362505365467
$ws.Range("AL40").Value = 0.6922790551338126
$ws.Range("AM40").Value = 0.96313217895925529
$ws.Range("AV40").Value = 0.70489334642046875
$ws.Range("AQ41").Value = 0.9423664764685169
$ws.Range("AN42").Value = 0.91125815282778955
$ws.Range("AQ42").Value = 0.66017088035908833
$ws.Range("AR43").Value = 0.77362023078426057
$ws.Range("AT44").Value = 0.78608753139860399
$ws.Range("AQ45").Value = 0.97998288299629377
$ws.Range("AU45").Value = 0.96919156166232356
$ws.Range("AS46").Value = 0.81960052671513517
$ws.Range("AU46").Value = 0.82778686176752558
$ws.Range("BI46").Value = 0.73273670042584182
$ws.Range("AV47").Value = 0.87107774640236191
$ws.Range("AW48").Value = 0.85277742659876599
$ws.Range("AX48").Value = 0.87064447785932741
$ws.Range("AY49").Value = 0.80786236865529126
$ws.Range("AW50").Value = 0.66661335312292958
$ws.Range("AX51").Value = 0.93821350983823382
$ws.Range("AZ51").Value = 0.83896299086194848
$ws.Range("BA51").Value = 0.99090491246333601
$ws.Range("E52").Value = 0.67226311434990327
$ws.Range("AW52").Value = 0.68956312151425925
$ws.Range("BA52").Value = 0.85230387931807416
$ws.Range("AM53").Value = 0.64355834072790696
$ws.Range("BA54").Value = 0.77194235520576515
$ws.Range("BD54").Value = 0.98687988503153545
$ws.Range("P55").Value = 0.57705022345777679
$ws.Range("BB55").Value = 0.96565176512634332
$ws.Range("BC57").Value = 0.95719864446298231
$ws.Range("BD57").Value = 0.7527139223971604
$ws.Range("F58").Value = 0.76486905041460962
$ws.Range("O58").Value = 0.81080660879069222
$ws.Range("Q58").Value = 0.68759165736200245
$ws.Range("BH58").Value = 0.91804064750767567
$ws.Range("BE59").Value = 0.92634616538044878
$ws.Range("BF59").Value = 0.83968561787498175
$ws.Range("BH59").Value = 0.60563917712368942
$ws.Range("BI60").Value = 0.99890973758808554
$ws.Range("BC61").Value = 0.99423546421850717
$ws.Range("BJ61").Value = 0.77620997589916252
$ws.Range("BH62").Value = 0.60446403277330207
$ws.Range("BK62").Value = 0.67329492325607676
$ws.Range("BL62").Value = 0.73226749806828928
$ws.Range("BB63").Value = 0.90221197869527781
$ws.Range("BI63").Value = 0.54377527441007301
$ws.Range("BM64").Value = 0.59901022939820625
$ws.Range("BN64").Value = 0.96054995255122977
$ws.Range("AW65").Value = 0.83535562461133295
$ws.Range("BO65").Value = 0.9922404187548115
$ws.Range("BP66").Value = 0.62409669713198546
$ws.Range("AR67").Value = 0.7442355695570726
$ws.Range("BN67").Value = 0.86770777673031407
$ws.Range("AT68").Value = 0.94457163427058832
$ws.Range("BO68").Value = 0.92463525661813195
